# Update InterFace imposible for VanLang Logo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the attendance date from July 19 to July 21, 2018
$ws.Range("C7").Formula = "=DATE(2018,7,21)"

# Update the time-in (F8) and time-out (G8) values
# F8: 9:00 AM -> 7:00 AM
$ws.Range("F8").Value = 7/24
# G8: 12:00 PM -> 9:30 PM
$ws.Range("G8").Value = 21.5/24

# Update the active selection on the sheet to F8
$ws.Range("F8").Select()
